$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 1.49
$ws.Range("H4").Value = 1.62
$ws.Range("K4").Value = 5.1
$ws.Range("G5").Value = 2.48
$ws.Range("S5").Value = 5.5
$ws.Range("F6").Value = 1.41
$ws.Range("G6").Value = 1.52
$ws.Range("H6").Value = 8.6
$ws.Range("I6").Value = 13.5
$ws.Range("J6").Value = 3.8
$ws.Range("K6").Value = 5.5
$ws.Range("M6").Value = 1.09
$ws.Range("N6").Value = 2.9
$ws.Range("O6").Value = 1.41
$ws.Range("P6").Value = 1.67
$ws.Range("Q6").Value = 2.04
$ws.Range("R6").Value = 1.25
$ws.Range("S6").Value = 3.75
$ws.Range("T6").Value = 2.44
$ws.Range("U6").Value = 1.56
$ws.Range("W6").Value = 2.86
$ws.Range("AN7").Value = 1000
$ws.Range("G7").Value = 1.11
$ws.Range("H7").Value = 26
$ws.Range("J7").Value = 7.8
$ws.Range("N7").Value = 3.5
$ws.Range("P7").Value = 3.5
$ws.Range("Q7").Value = 1.27
$ws.Range("R7").Value = 2.08
$ws.Range("S7").Value = 1.54
$ws.Range("U7").Value = 1.04
$ws.Range("W7").Value = 8.6
$ws.Range("F8").Value = 1.32
$ws.Range("I8").Value = 13
$ws.Range("J8").Value = 4.6
$ws.Range("K8").Value = 8
$ws.Range("P8").Value = 2.86
$ws.Range("U8").Value = 2.02
$ws.Range("V8").Value = 1.09
$ws.Range("AC9").Value = 990
$ws.Range("F9").Value = 3.15
$ws.Range("I9").Value = 2.44
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 3.95
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 3.85
$ws.Range("O9").Value = 1.28
$ws.Range("P9").Value = 1.98
$ws.Range("Q9").Value = 1.84
$ws.Range("R9").Value = 1.38
$ws.Range("S9").Value = 3.1
$ws.Range("T9").Value = 1.69
$ws.Range("U9").Value = 2.2
$ws.Range("X9").Value = 19
$ws.Range("T10").Value = 1.66
$ws.Range("T11").Value = 1.53
$ws.Range("F12").Value = 2.32
$ws.Range("I12").Value = 3.35
$ws.Range("N12").Value = 3.15
$ws.Range("O12").Value = 1.37
$ws.Range("R12").Value = 1.27
$ws.Range("T12").Value = 1.81
$ws.Range("U12").Value = 1.96
$ws.Range("V12").Value = 1.42
$ws.Range("AC14").Value = 9
$ws.Range("F14").Value = 2.14
$ws.Range("N14").Value = 3.9
$ws.Range("O14").Value = 1.27
$ws.Range("P14").Value = 1.99
$ws.Range("Q14").Value = 1.82
$ws.Range("R14").Value = 1.37
$ws.Range("T14").Value = 1.71
$ws.Range("U14").Value = 2.18
$ws.Range("Y14").Value = 990
$ws.Range("N15").Value = 1.1
$ws.Range("P15").Value = 3.15
$ws.Range("R15").Value = 1.76
$ws.Range("AB16").Value = 13.5
$ws.Range("AC16").Value = 10.5
$ws.Range("AD16").Value = 18.5
$ws.Range("AF16").Value = 18
$ws.Range("AG16").Value = 13.5
$ws.Range("AH16").Value = 20
$ws.Range("AI16").Value = 980
$ws.Range("AJ16").Value = 34
$ws.Range("AK16").Value = 26
$ws.Range("AN16").Value = 17.5
$ws.Range("F16").Value = 2.16
$ws.Range("G16").Value = 2.38
$ws.Range("H16").Value = 3.15
$ws.Range("I16").Value = 3.6
$ws.Range("J16").Value = 3.5
$ws.Range("K16").Value = 4
$ws.Range("O16").Value = 1.25
$ws.Range("T16").Value = 1.64
$ws.Range("U16").Value = 2.26
$ws.Range("V16").Value = 1.38
$ws.Range("W16").Value = 1.73
$ws.Range("Z16").Value = 34
$ws.Range("J17").Value = 3.85
$ws.Range("AK18").Value = 36
$ws.Range("L18").Value = 1.43
$ws.Range("O18").Value = 1.35
$ws.Range("AC20").Value = 16.5
$ws.Range("AN20").Value = 50
$ws.Range("AO20").Value = 5.1
$ws.Range("G20").Value = 6.8
$ws.Range("I20").Value = 1.58
$ws.Range("J20").Value = 4.9
$ws.Range("K20").Value = 5.9
$ws.Range("L20").Value = 1.15
$ws.Range("M20").Value = 1.02
$ws.Range("N20").Value = 7
$ws.Range("O20").Value = 1.13
$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 1.37
$ws.Range("R20").Value = 1.82
$ws.Range("S20").Value = 1.98
$ws.Range("T20").Value = 1.57
$ws.Range("U20").Value = 2.4
$ws.Range("V20").Value = 2.72
$ws.Range("X20").Value = 42
$ws.Range("Y20").Value = 17.5
$ws.Range("Z20").Value = 15.5
$ws.Range("G21").Value = 5.9
$ws.Range("Q22").Value = 2.46
$ws.Range("F23").Value = 2.14
$ws.Range("K23").Value = 5.5
$ws.Range("N23").Value = 1.1
$ws.Range("T23").Value = 1.55
$ws.Range("AA24").Value = 1000
$ws.Range("AB24").Value = 18
$ws.Range("AF24").Value = 12.5
$ws.Range("AH24").Value = 32
$ws.Range("AJ24").Value = 14.5
$ws.Range("H24").Value = 14
$ws.Range("J24").Value = 7.8
$ws.Range("K24").Value = 9.4
$ws.Range("P24").Value = 3.8
$ws.Range("Q24").Value = 1.29
$ws.Range("S24").Value = 1.76
$ws.Range("T24").Value = 1.67
$ws.Range("U24").Value = 2.16
$ws.Range("Y24").Value = 990
$ws.Range("AB25").Value = 15
$ws.Range("J25").Value = 3.7
$ws.Range("K25").Value = 3.75
$ws.Range("N25").Value = 5.1
$ws.Range("O25").Value = 1.23
$ws.Range("P25").Value = 2.42
$ws.Range("S25").Value = 2.68
